# update_count function and the section where the "All SANs" sheet is appended.
# These changes ensure that the log is written in the correct order and that
# the log view is refreshed immediately after each operation, displaying the
# entries sorted by the timestamp in descending order ("newest first").
# Additionally, this corrects the order of data appended to the "All SANs" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "4.2 Items" - Desktop Mini G9 count bumped from 90/91 to 93/94
# ---------------------------------------------------------------------------
$items42 = $wb.Worksheets.Item("4.2 Items")
$items42.Range("B2").Value = 93
$items42.Range("C2").Value = 94

# ---------------------------------------------------------------------------
# 2) "4.2 Timestamps" - three new "add" log rows for Desktop Mini G9
# ---------------------------------------------------------------------------
$ts42 = $wb.Worksheets.Item("4.2 Timestamps")
$ts42.Range("A17").Value = "Desktop Mini G9"
$ts42.Range("B17").Value = "add"
$ts42.Range("C17").Value = "SAN888444"
$ts42.Range("D17").Value = "2023-12-28 23:20:29"

$ts42.Range("A18").Value = "Desktop Mini G9"
$ts42.Range("B18").Value = "add"
$ts42.Range("C18").Value = "SAN555555"
$ts42.Range("D18").Value = "2023-12-28 23:21:36"

$ts42.Range("A19").Value = "Desktop Mini G9"
$ts42.Range("B19").Value = "add"
$ts42.Range("C19").Value = "SAN344556"
$ts42.Range("D19").Value = "2023-12-28 23:32:03"

# ---------------------------------------------------------------------------
# 3) "BR Items" - Desktop Mini G9 count bumped from 8/9 to 9/12
# ---------------------------------------------------------------------------
$itemsBR = $wb.Worksheets.Item("BR Items")
$itemsBR.Range("B2").Value = 9
$itemsBR.Range("C2").Value = 12

# ---------------------------------------------------------------------------
# 4) "BR Timestamps" - three new "add" log rows for Desktop Mini G9
# ---------------------------------------------------------------------------
$tsBR = $wb.Worksheets.Item("BR Timestamps")
$tsBR.Range("A2").Value = "Desktop Mini G9"
$tsBR.Range("B2").Value = "add"
$tsBR.Range("C2").Value = "SAN343434"
$tsBR.Range("D2").Value = "2023-12-29 00:07:07"

$tsBR.Range("A3").Value = "Desktop Mini G9"
$tsBR.Range("B3").Value = "add"
$tsBR.Range("C3").Value = "SAN454545"
$tsBR.Range("D3").Value = "2023-12-29 00:07:10"

$tsBR.Range("A4").Value = "Desktop Mini G9"
$tsBR.Range("B4").Value = "add"
$tsBR.Range("C4").Value = "SAN535353"
$tsBR.Range("D4").Value = "2023-12-29 00:07:13"

# ---------------------------------------------------------------------------
# 5) "All SANs" - corrected order data appended (SAN Number, Item, Time)
# ---------------------------------------------------------------------------
$allSans = $wb.Worksheets.Item("All SANs")

$allSans.Range("A17").Value = "Desktop Mini G9"
$allSans.Range("B17").Value = "SAN888444"
$allSans.Range("C17").Value = "2023-12-28 23:20:29"

$allSans.Range("A18").Value = "Desktop Mini G9"
$allSans.Range("B18").Value = "SAN555555"
$allSans.Range("C18").Value = "2023-12-28 23:21:36"

$allSans.Range("A19").Value = "Desktop Mini G9"
$allSans.Range("B19").Value = "SAN344556"
$allSans.Range("C19").Value = "2023-12-28 23:32:03"

$allSans.Range("A20").Value = "Desktop Mini G9"
$allSans.Range("B20").Value = "SAN343434"
$allSans.Range("C20").Value = "2023-12-29 00:07:07"

$allSans.Range("A21").Value = "Desktop Mini G9"
$allSans.Range("B21").Value = "SAN454545"
$allSans.Range("C21").Value = "2023-12-29 00:07:10"

$allSans.Range("A22").Value = "Desktop Mini G9"
$allSans.Range("B22").Value = "SAN535353"
$allSans.Range("C22").Value = "2023-12-29 00:07:13"
